$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# New order rows 42-51, appended after the existing data (row 41 was the last one).
# Column A = PackageID (only set on rows that start a new package),
# Column C = FlowerName, Column F = Number.
# Format the numeric-looking columns (A, F) as text first so values like "5"
# or "7" are stored as text, matching every other data cell on this sheet.
$ws.Range("A42:A51").NumberFormat = "@"
$ws.Range("F42:F51").NumberFormat = "@"

$ws.Range("C42").Value = "683_锦鲤红_undefined_undefined_1bunch"
$ws.Range("F42").Value = "5"

$ws.Range("C43").Value = "681_锦鲤橙_undefined_undefined_1bunch"
$ws.Range("F43").Value = "5"

$ws.Range("A44").Value = "7"
$ws.Range("C44").Value = "546_非洲菊绿宝石_undefined_undefined_1bunch"
$ws.Range("F44").Value = "5"

$ws.Range("C45").Value = "612_康乃馨古巴爱情_cuba love_undefined_20stems"
$ws.Range("F45").Value = "5"

$ws.Range("C46").Value = "386_菟葵绿粉 `ngreen_undefined_undefined_1bunch"
$ws.Range("F46").Value = "8"

$ws.Range("C47").Value = "48_香格里拉_undefined_Gerbera L._10stems"
$ws.Range("F47").Value = "10"

$ws.Range("C48").Value = "611_康乃馨奶油白_cream white_undefined_20stems"
$ws.Range("F48").Value = "5"

$ws.Range("C49").Value = "43_拉丝红_Spider Red_Gerbera L._20stems"
$ws.Range("F49").Value = "10"

$ws.Range("C50").Value = "574_迷你菊白_undefined_undefined_1bunch"
$ws.Range("F50").Value = "10"

$ws.Range("C51").Value = "70_朝霞mini_undefined_Gerbera L._20stems"

# Summary sheet: TotalNumber string is the concatenation of all Number values;
# recomputed now that the Orders sheet has 10 more rows of data.
$summary.Range("G2").NumberFormat = "@"
$summary.Range("G2").Value = "0146137101398786310151510155101051010155510755555101015151010105555810510100"
